# Interdiff between v5 and v6:
# The cached text of every "datetimeFigureOut" Date placeholder field
# (slide master, every slide layout, and the notes master) is refreshed
# from 4/4/2018 to 4/16/2018.

$p = $ppt.ActivePresentation
$newDate = "4/16/2018"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }
        if ($isDatePlaceholder) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout attached to the slide master
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Notes master
if ($p.HasNotesMaster) {
    $notesMaster = $p.NotesMaster
    Update-DatePlaceholder $notesMaster.Shapes
}
